$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(newPrice, newVolume)  ($null means "unchanged, skip")
$updates = @{
    2  = @("27.908.27", "  +0.92%  ")
    3  = @("1.886.56", "  +0.27%  ")
    4  = @($null, "  +1.77%  ")
    5  = @("335.85", "  +1.41%  ")
    6  = @("1.019", $null)
    7  = @("0.4652", "  -1.67%  ")
    8  = @("0.3908", "  -1.47%  ")
    9  = @("47.06", "  -0.44%  ")
    10 = @("0.07953", "  -1.05%  ")
    11 = @("1.012", "  -1.16%  ")
    12 = @("21.58", "  -0.89%  ")
    13 = @("1.890.65", "  +1.60%  ")
    14 = @("5.921", "  -0.80%  ")
    15 = @("7.085", "  -1.20%  ")
    16 = @($null, "  +1.93%  ")
    17 = @("0.06763", "  +2.08%  ")
    18 = @("86.88", "  -0.07%  ")
    19 = @("0.00001044", "  +0.10%  ")
    20 = @("16.96", "  -1.32%  ")
    21 = @("1.018", "  +1.66%  ")
    22 = @("27.933.92", "  +0.94%  ")
    23 = @("5.481", "  -0.40%  ")
    24 = @("10.91", "  -0.58%  ")
    25 = @("2.350", "  +1.78%  ")
    26 = @("2.112.52", "  +1.47%  ")
    27 = @("159.06", "  +2.08%  ")
    28 = @($null, "  -1.33%  ")
    29 = @("2.066", "  -1.30%  ")
    30 = @("5.409", "  -2.82%  ")
    31 = @("121.20", "  -0.85%  ")
    32 = @("0.9594", "  -0.80%  ")
    33 = @("0.09481", "  -0.64%  ")
    34 = @("3.674", "  +1.15%  ")
    35 = @("1.372", "  -6.15%  ")
    36 = @("5.324", "  +0.26%  ")
    37 = @("0.06086", "  -0.62%  ")
    38 = @($null, "  -0.78%  ")
    39 = @("1.216", "  -1.31%  ")
    40 = @("8.062", "  -0.97%  ")
    41 = @("0.5936", "  -1.15%  ")
    42 = @("0.1881", "  -0.95%  ")
    43 = @("10.22", "  -0.74%  ")
    44 = @("1.273", "  +1.91%  ")
    45 = @("0.5650", "  -1.13%  ")
    46 = @("12.16", "  -0.78%  ")
    47 = @("3.393", "  -0.40%  ")
    48 = @("1.915", "  -0.93%  ")
    49 = @("0.06906", "  +1.27%  ")
    50 = @("113.55", "  +2.45%  ")
    51 = @("1.064", "  -0.39%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($null -ne $priceVal) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $priceVal
        $cell.Style = "Normal"
    }

    if ($null -ne $volVal) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $volVal
        $cell.Style = "Normal"
    }
}
